$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 551
$ws1.Range("F4").Value = 181
$ws1.Range("F5").Value = 300
$ws1.Range("F8").Value = 2338
$ws1.Range("F10").Value = 5864
$ws1.Range("F12").Value = 381
$ws1.Range("F13").Value = 6

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 15

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 551
$ws4.Range("F5").Value = 181
$ws4.Range("F6").Value = 300
$ws4.Range("F10").Value = 15
$ws4.Range("F11").Value = 2338
$ws4.Range("F13").Value = 5864
$ws4.Range("F15").Value = 381
$ws4.Range("F17").Value = 6
